$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet: insert a new blank column before column N.
# This shifts the existing "Late" / "Heading" / "Outstanding" columns
# (previously N/O/P) one slot to the right, becoming O/P/Q, and widens
# the used range from A1:P14 to A1:Q14.
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N").Insert()

# The workbook now opens with "Repayment Schedule" as the active tab
# (rather than "Input"), scrolled/selected at M15.
[void]$wsSchedule.Activate()
[void]$wsSchedule.Range("M15").Select()
